$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Turn the first paragraph ("Ejemplo de estructura: ") into a new
#    title paragraph ("Estructura interna de un repositorio"),
#    pushing the original text down into a newly inserted paragraph
#    that keeps the original formatting. Also move the "_GoBack"
#    bookmark so it sits right after the new title text.
# ------------------------------------------------------------------

$firstPara = $d.Paragraphs(1)

# Insert a brand-new paragraph right before the current first one;
# it inherits the paragraph/run formatting (pPr + rPr) automatically.
$firstPara.Range.InsertParagraphBefore() | Out-Null

# The former paragraph 1 ("Ejemplo de estructura: ") is now paragraph 2.
# Paragraph 1 is empty (just inherited formatting) - give it the new title.
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Text = "Estructura interna de un repositorio"

# Re-fetch the (now retitled) paragraph and append a throw-away marker
# character so we can drop a zero-width bookmark right after the real
# text without it being absorbed into the following paragraph - then
# remove the marker again.
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.InsertAfter("@")
$titlePara = $d.Paragraphs(1)
$titleEnd = $titlePara.Range.End
$bmRange = $d.Range($titleEnd - 2, $titleEnd - 2)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
$markerRange = $d.Range($titleEnd - 2, $titleEnd - 1)
$markerRange.Delete() | Out-Null

# ------------------------------------------------------------------
# 2) Drop the stale <w:lastRenderedPageBreak/> cached before the
#    "OBJETIVO/S DE LA ACTIVIDAD ..." run, without touching its text.
# ------------------------------------------------------------------

$targetText = "OBJETIVO/S DE LA ACTIVIDAD (Relacionados con la unidad/m"
$found = $d.Content.Find.Execute($targetText, $true, $false, $false, $false, $false, $true, 1, $false, $targetText, 2)
if (-not $found) {
    throw "Could not locate the OBJETIVO/S paragraph to clean up the page break marker"
}

# ------------------------------------------------------------------
# 3) The former last paragraph only held the stray "_GoBack" bookmark;
#    it has already been relocated above, so that paragraph is left
#    empty automatically. Nothing further required here.
# ------------------------------------------------------------------

Write-Output "done"
